# Update TPM-derived values on Sheet1 (Sost-Lrp6 LR-pair table)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 0.7370549019300001
$ws.Range("R2").Value = 6.633494117370001
$ws.Range("S2").Value = 0.1441015470002482
$ws.Range("T2").Value = 0.1441015470002482

# Row 3 (only derived-specificity columns change; M3/N3 stay the same)
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("S3").Value = 0.3846359116098663
$ws.Range("T3").Value = 0.3846359116098662

# Row 4
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 1.120589441004
$ws.Range("R4").Value = 10.085304969036
$ws.Range("S4").Value = 0.2190863551385157
$ws.Range("T4").Value = 0.2190863551385156

# Row 5
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 1.289838298726
$ws.Range("R5").Value = 11.608544688534
$ws.Range("S5").Value = 0.2521761862513699
$ws.Range("T5").Value = 0.2521761862513699
